$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift "hour" column (B) values by +1, starting at 9 instead of 8 (rows 2-16)
$ws.Range("B2").Value = 9
$ws.Range("B3").Value = 10
$ws.Range("B4").Value = 11
$ws.Range("B5").Value = 12
$ws.Range("B6").Value = 13
$ws.Range("B7").Value = 14
$ws.Range("B8").Value = 15
$ws.Range("B9").Value = 16
$ws.Range("B10").Value = 17
$ws.Range("B11").Value = 18
$ws.Range("B12").Value = 19
$ws.Range("B13").Value = 20
$ws.Range("B14").Value = 21
$ws.Range("B15").Value = 22
$ws.Range("B16").Value = 23

# Update "Offerte" column (C) with new pasthour record values (only updated at 00, rest cleared)
$ws.Range("C2").Value = 52
$ws.Range("C3").Value = 103
$ws.Range("C4").Value = 140
$ws.Range("C5").Value = 130
$ws.Range("C6").Value = 128
$ws.Range("C7").Value = 111
$ws.Range("C8").Value = 90
$ws.Range("C9").Value = 106
$ws.Range("C10").Value = 109
